$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... Allow reuse of this data across applications by sharing these
#    files." -> "... by storing a copy in a common location."
# ---------------------------------------------------------------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    "sharing these files",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "storing a copy in a common location",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Change the glyph that precedes "Addition of Words animated to
#    reinforce what's going on." from the "Arial Unicode MS" circled-2
#    symbol to a "MS Mincho" check mark.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Addition of Words animated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $anchor.Paragraphs(1).Range.Start
$symbolRange = $d.Range($paraStart, $paraStart)

$symbolRange.Find.ClearFormatting()
$symbolRange.Find.Replacement.ClearFormatting()
$symbolRange.Find.Replacement.Font.Name = "MS Mincho"
$symbolRange.Find.Replacement.Font.NameFarEast = "MS Mincho"
$symbolRange.Find.Replacement.Font.NameBi = "MS Mincho"
$symbolRange.Find.Execute(
    "➁",
    $true, $false, $false, $false, $false,
    $true, 0, $true,
    "✔",
    2
) | Out-Null

# Reset the Find/Replacement formatting criteria so later plain-text
# replacements below do not inherit the MS Mincho formatting just used.
$symbolRange.Find.ClearFormatting()
$symbolRange.Find.Replacement.ClearFormatting()

# ---------------------------------------------------------------------------
# 3) "Print draft dictionary. (In Sept 2007 doesn't output pictures or
#    custom sorting)" -> "Print draft dictionary. (In Oct 2007 doesn't
#    output pictures)"
# ---------------------------------------------------------------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    "Sept",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Oct",
    2
) | Out-Null

$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    " or custom sorting)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ")",
    2
) | Out-Null
